$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.144.02"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "2.304.32"
$ws.Range("E3").Value = "  +4.19%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "252.89"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "0.641"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("D7").Value = "74.48"
$ws.Range("E7").Value = "  +10.02%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "0.649"
$ws.Range("E9").Value = "  +5.46%  "
$ws.Range("D10").Value = "39.74"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").Value = "0.0991"
$ws.Range("E11").Value = "  +5.76%  "
$ws.Range("D12").Value = "59.33"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "7.39"
$ws.Range("E13").Value = "  +5.24%  "
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "2.638.04"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("D16").Value = "15.52"
$ws.Range("E16").Value = "  +7.17%  "
$ws.Range("D17").Value = "0.884"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").Value = "2.278.13"
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").Value = "42.940.56"
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  +4.73%  "
$ws.Range("D21").Value = "6.33"
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").Value = "72.73"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "236.24"
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("E24").Value = "  +9.82%  "
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("D26").Value = "11.65"
$ws.Range("E26").Value = "  +3.69%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "2.44"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").Value = "167.44"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.20"
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("D33").Value = "6.44"
$ws.Range("E33").Value = "  +9.56%  "
$ws.Range("D34").Value = "0.128"
$ws.Range("E34").Value = "  +5.17%  "
$ws.Range("D35").Value = "0.0821"
$ws.Range("E35").Value = "  +5.35%  "
$ws.Range("D36").Value = "32.31"
$ws.Range("E36").Value = "  +25.05%  "
$ws.Range("E37").Value = "  +4.68%  "
$ws.Range("D38").Value = "4.74"
$ws.Range("E38").Value = "  +15.46%  "
$ws.Range("D39").Value = "4.79"
$ws.Range("E39").Value = "  +4.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0310"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.60"
$ws.Range("E41").Value = "  +21.43%  "
$ws.Range("D42").Value = "2.35"
$ws.Range("E42").Value = "  +5.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.00"
$ws.Range("E43").Value = "  +6.22%  "
$ws.Range("D44").Value = "0.215"
$ws.Range("E44").Value = "  +9.72%  "
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").Value = "62.38"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "9.15"
$ws.Range("E46").Value = "  +7.08%  "
$ws.Range("D47").Value = "4.88"
$ws.Range("E47").Value = "  -4.77%  "
$ws.Range("D48").Value = "0.104"
$ws.Range("E48").Value = "  +4.23%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").Value = "1.19"
$ws.Range("E50").Value = "  +3.17%  "
$ws.Range("D51").Value = "98.52"
$ws.Range("E51").Value = "  +6.12%  "
